# Yarpiz PSO adaptation: rework the DE results table to PSO's MaxFES grid
# and drop the 51st run ("Run 50"), recomputing the Mean column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Header: "Gen" -> "MaxFES"
$ws.Cells.Item(1, 1).Value = "MaxFES"

# 2) Column A (rows 2..14): generation counts -> normalized MaxFES fractions
$maxfes = @(0, 0.001, 0.01, 0.1, 0.2, 0.3, 0.4, 0.5, 0.6, 0.7, 0.8, 0.9, 1)
for ($i = 0; $i -lt $maxfes.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $maxfes[$i]
}

# 3) Drop the "Run 50" column (AZ). This shifts the trailing "Mean" column
#    (previously BA) left into AZ, and the sheet's used range shrinks from
#    A1:BA14 to A1:AZ14.
$ws.Columns("AZ").Delete()

# 4) Recompute the Mean column (now AZ) for rows 2..14, excluding the
#    removed 51st run.
$meanVals = @(
    273.41293236,
    212.72956712,
    25.8243338,
    0.35793385,
    0.16205779,
    0.11615425,
    0.09050233000000001,
    0.07088042999999999,
    0.06289893000000001,
    0.05429183,
    0.05008235,
    0.04649491,
    0.04365521
)
for ($i = 0; $i -lt $meanVals.Length; $i++) {
    $ws.Cells.Item($i + 2, 52).Value = $meanVals[$i]
}
